$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: new "inflation" time-variant data row ---

# A37: label cell. Style mirrors the other label cells (A33:A36, style
# index 16 = bold font + border) but WITHOUT the blue fill, i.e. style 21.
$ws.Range("A37").Value = "inflation"
$ws.Range("A37").Font.Bold = $true
$ws.Range("A37").Interior.Pattern = 0
$ws.Range("A37").Borders(10).LineStyle = 1

# B37:D37 - flag columns, plain values; these pick up the column default
# style (index 2) automatically, same as existing rows.
$ws.Range("B37").Value = "yes"
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = "yes"

# L37:M37 - first two data points; these pick up the column default style
# (index 1) automatically, same as the header/column formatting.
$ws.Range("L37").Value = 0.041
$ws.Range("M37").Value = 0.131

# N37:AZ37 - bulk of the yearly data, entered with no explicit formatting
# (left at the workbook default / "Normal" style so no style index is
# written for these cells).
$dataRange = $ws.Range("N37:AZ37")
$dataRange.Value = 0.05
$dataRange.Style = "Normal"

$values = @{
    "N37"  = 0.114
    "O37"  = 0.07
    "P37"  = 0.061
    "Q37"  = 0.078
    "R37"  = 0.145
    "S37"  = 0.112
    "T37"  = 0.07
    "U37"  = 0.067
    "V37"  = 0.053
    "W37"  = 0.044
    "X37"  = 0.018
    "Y37"  = 0.057
    "Z37"  = 0.118
    "AA37" = 0.062
    "AB37" = 0.082
    "AC37" = 0.065
    "AD37" = 0.049
    "AE37" = 0.052
    "AF37" = 0.008
    "AG37" = 0.022
    "AH37" = 0.031
    "AI37" = 0.034
    "AJ37" = 0.057
    "AK37" = 0.02
    "AL37" = 0.011
    "AM37" = 0.043
    "AN37" = 0.008
    "AO37" = 0.042
    "AP37" = 0.028
    "AQ37" = 0.024
    "AR37" = 0.025
    "AS37" = 0.048
    "AT37" = 0.077
    "AU37" = 0.032
    "AV37" = 0.037
    "AW37" = 0.073
    "AX37" = 0.034
    "AY37" = 0.029
    "AZ37" = 0.005
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

# Re-apply the "no explicit style" formatting after writing the values
# above (writing .Value again can re-stamp the column default style).
$dataRange.Style = "Normal"

# BA37 - final data point, styled with the highlighted style used
# elsewhere in the sheet (index 7). Copy format from an existing cell
# using that same style so the original style index is reused instead of
# a new one being minted.
$ws.Range("BA37").Value = 0.014
$ws.Range("J2").Copy() | Out-Null
$ws.Range("BA37").PasteSpecial(-4122) | Out-Null

# --- Data validation: extend the first decimal validation to cover the
# new row's N37:AW37 range (matches the diff's sqref addition). ---
$dv = $ws.Range("E28:BA31 E2:BC2 E4:BC4 E10:BC10 E20:BC20 E22:BC22 E24:BC24 E26:BC26 E12:BC12 E14:BC14 E16:BC16 E18:BC18 E8:BC8 BB28:BC30")
$dv.Validation.Delete()
$dv.Validation.Add(5, 1, 1, 0, 100)
$dv2 = $ws.Range("N37:AW37")
$dv2.Validation.Delete()
$dv2.Validation.Add(5, 1, 1, 0, 100)

# --- Update the active selection (bottomRight pane) to J37 to match the
# saved workbook view state. ---
[void]$ws.Range("J37").Select()
